$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "1SF0"
$ws.Range("B2").Value = "SDMHF"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = "Femenino"
$ws.Range("E2").Value = "18/2/2021 - 07:42:06 p. m."
$ws.Range("F2").Value = "KN"
